# Insert a new weekly record as row 15 in the "Achicoria" price sheet.
# All rows from the old row 15 through the old row 62 shift down by one
# (to rows 16..63) and the newly inserted row 15 receives the latest
# week's data (same market/region/category metadata, new date & price
# figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 15 (and everything below it) down by one row.
$ws.Rows.Item(15).Insert()

# Populate the newly-inserted row 15 with the new weekly record.
$ws.Range("A15").Value = 9
$ws.Range("B15").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C15").Value = "Metropolitana"
$ws.Range("D15").Value = 45071
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 100112010
$ws.Range("G15").Value = "Achicoria"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 90
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("N15").Value = "$/caja 16 unidades"
$ws.Range("O15").Value = "Provincia de Quillota"
$ws.Range("P15").Value = 438
$ws.Range("Q15").Value = 16
$ws.Range("R15").Value = "Hortaliza"
